# Update the build timestamp embedded in the version string throughout
# the workbook, from "February 03 2026 17.29.55 EST" to
# "February 03 2026 18.05.36 EST".

$wb = $excel.ActiveWorkbook

$newStamp = "February 03 2026 18.05.36 EST"

# --- "About" sheet -------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

$wsAbout.Range("A6").Value = "Recommended Citation:  " + '"' + "Global Energy Monitor, Coal mine boundaries and methane sources for Flying Eagle Mine, United States, M3466, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)" + '"'

# --- "Boundaries and methane sources" sheet -------------------------------
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 7; $row++) {
    $cell = $wsData.Range("S$row")
    $cell.Value = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"
}
